$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert DJ600 (beta 0.75) above EMB (original row 9)
$ws.Rows.Item(9).Insert()
$ws.Range("A9").Value = "DJ600"
$ws.Range("B9").Value = 0.75

# Insert GLEN (beta 0.38) above GDXJ (original row 18, now shifted to row 19)
$ws.Rows.Item(19).Insert()
$ws.Range("A19").Value = "GLEN"
$ws.Range("B19").Value = 0.38

# Insert PAAS (beta 0.123, comment "basically SLV") above RACE (original row 27, now shifted to row 29)
$ws.Rows.Item(29).Insert()
$ws.Range("A29").Value = "PAAS"
$ws.Range("B29").Value = 0.123
$ws.Range("D29").Value = "basically SLV"

# Insert SLV (value = SI's beta, via formula) above STNG (original row 34, now shifted to row 37)
$ws.Rows.Item(37).Insert()
$ws.Range("A37").Value = "SLV"
$ws.Range("B37").Formula = "=B36"

# Fix typo in SI's comment: "depended" -> "dependent"
$ws.Range("D36").Value = "unicorn bay, but very term dependent. back to 2014 it drops to .1"

# Update the active cell / selection to match the author's final cursor position
$ws.Range("I32").Select() | Out-Null
